# Handle "no suspicious line detected" cases: populate the results table
# with a bold/bordered/centered header row and two data rows.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$headers = @("Project", "Submission", "Patch", "Passed Tests")
$cols = @("A", "B", "C", "D")

for ($i = 0; $i -lt 4; $i++) {
    $ws.Range($cols[$i] + "1").Value = $headers[$i]
}

# Header row formatting: bold font, thin border all around, centered
# horizontally and top-aligned vertically. Format A1 first, then copy its
# resulting format onto the rest of the header row so every header cell
# shares the exact same style entry.
$a1 = $ws.Range("A1")
$a1.Font.Bold = $true
$a1.HorizontalAlignment = -4108  # xlCenter
$a1.VerticalAlignment = -4160    # xlTop
$a1.Borders.LineStyle = 1        # xlContinuous
$a1.Borders.Weight = 2           # xlThin

$a1.Copy()
$ws.Range("B1:D1").PasteSpecial(-4122)  # xlPasteFormats

# Data row 2 - "Submission" stored as a genuine number (59).
$ws.Range("A2").Value = "Cafe"
$ws.Range("B2").Value = 59
$ws.Range("C2").Value = "AstorMain-Cafe-59-1"
$ws.Range("D2").Value = "57 / 67"

# Data row 3 - "Submission" stored as text "59" (not a number) this time.
$ws.Range("A3").Value = "Cafe"

$scratch = $ws.Range("Z100")
$scratch.NumberFormat = "@"   # force text interpretation
$scratch.Value = "59"
$scratch.Copy()
$ws.Range("B3").PasteSpecial(-4163)  # xlPasteValues (keeps B3 unstyled)
$scratch.Clear()

$ws.Range("C3").Value = "AstorMain-Cafe-59-2"
$ws.Range("D3").Value = "58 / 67"
